$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 535.97437
$ws.Range("J17").Value = 287.19354
$ws.Range("L17").Value = 861.58062
$ws.Range("N17").Value = -1197.58062
$ws.Range("H18").Value = 3866.3333
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 3866.3333
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 3866.3333
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -4434.3333
$ws.Range("H26").Value = 44999.668
$ws.Range("J26").Value = 44999.668
$ws.Range("L26").Value = 44999.668
$ws.Range("N26").Value = -45687.668
$ws.Range("H51").Value = 6532.3335
$ws.Range("I51").Value = 2299
$ws.Range("J51").Value = 8649
$ws.Range("K51").Value = 2299
$ws.Range("L51").Value = 8649
$ws.Range("M51").Value = -1815
$ws.Range("N51").Value = -9617
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H127").Value = 1348.3846
$ws.Range("I127").Value = 577.8570999999999
$ws.Range("J127").Value = 2247.3333
$ws.Range("K127").Value = 1733.5713
$ws.Range("L127").Value = 6741.999899999999
$ws.Range("M127").Value = 3226.4287
$ws.Range("N127").Value = -16661.9999
$ws.Range("H129").Value = 1292.9747
$ws.Range("I129").Value = 469.4
$ws.Range("J129").Value = 1348.6216
$ws.Range("K129").Value = 1408.2
$ws.Range("L129").Value = 4045.8648
$ws.Range("M129").Value = 3591.8
$ws.Range("N129").Value = -14045.8648
$ws.Range("H132").Value = 22533084
$ws.Range("I132").Value = 24487078
$ws.Range("K132").Value = 73461234
$ws.Range("M132").Value = -73458704
$ws.Range("H137").Value = 673612.6
$ws.Range("I137").Value = 1908128.5
$ws.Range("J137").Value = 2680.152
$ws.Range("K137").Value = 5724385.5
$ws.Range("L137").Value = 8040.456
$ws.Range("M137").Value = -5721835.5
$ws.Range("N137").Value = -13140.456
$ws.Range("H138").Value = 2529.6572
$ws.Range("J138").Value = 3610.7778
$ws.Range("L138").Value = 10832.3334
$ws.Range("N138").Value = -21112.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 2647.6365
$ws.Range("I122").Value = 903.8333
$ws.Range("K122").Value = 2711.4999
$ws.Range("M122").Value = -261.4998999999998
$ws.Range("H132").Value = 2004.6364
$ws.Range("I132").Value = 1553.5435
$ws.Range("J132").Value = 4310.222
$ws.Range("K132").Value = 4660.6305
$ws.Range("L132").Value = 12930.666
$ws.Range("M132").Value = -2130.6305
$ws.Range("N132").Value = -17990.666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H23").Value = 32503
$ws.Range("I23").Value = 1012
$ws.Range("J23").Value = 43000
$ws.Range("K23").Value = 1012
$ws.Range("L23").Value = 43000
$ws.Range("M23").Value = -729
$ws.Range("N23").Value = -43566
$ws.Range("H86").Value = 2476.25
$ws.Range("I86").Value = 2520
$ws.Range("J86").Value = 2345
$ws.Range("K86").Value = 2520
$ws.Range("L86").Value = 2345
$ws.Range("M86").Value = -1397
$ws.Range("N86").Value = -4591
$ws.Range("H89").Value = 2476.25
$ws.Range("I89").Value = 2520
$ws.Range("J89").Value = 2345
$ws.Range("K89").Value = 12600
$ws.Range("L89").Value = 11725
$ws.Range("M89").Value = -6984
$ws.Range("N89").Value = -22957
$ws.Range("H134").Value = 2748.0925
$ws.Range("I134").Value = 1097.6207
$ws.Range("J134").Value = 4662.64
$ws.Range("K134").Value = 3292.8621
$ws.Range("L134").Value = 13987.92
$ws.Range("M134").Value = -757.8620999999998
$ws.Range("N134").Value = -19057.92

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 12928.429
$ws.Range("I23").Value = 6908.909
$ws.Range("J23").Value = 35000
$ws.Range("K23").Value = 6908.909
$ws.Range("L23").Value = 35000
$ws.Range("M23").Value = -6668.909
$ws.Range("N23").Value = -35480
$ws.Range("H27").Value = 12928.429
$ws.Range("I27").Value = 6908.909
$ws.Range("J27").Value = 35000
$ws.Range("K27").Value = 6908.909
$ws.Range("L27").Value = 35000
$ws.Range("M27").Value = -6716.909
$ws.Range("N27").Value = -35384
$ws.Range("H31").Value = 239475.02
$ws.Range("I31").Value = 501370.8
$ws.Range("J31").Value = 3768.8
$ws.Range("K31").Value = 501370.8
$ws.Range("L31").Value = 3768.8
$ws.Range("M31").Value = -501075.8
$ws.Range("N31").Value = -4358.8
$ws.Range("H34").Value = 239475.02
$ws.Range("I34").Value = 501370.8
$ws.Range("J34").Value = 3768.8
$ws.Range("K34").Value = 501370.8
$ws.Range("L34").Value = 3768.8
$ws.Range("M34").Value = -501168.8
$ws.Range("N34").Value = -4172.8
$ws.Range("H58").Value = 2455.0278
$ws.Range("I58").Value = 1299.3572
$ws.Range("J58").Value = 6499.875
$ws.Range("K58").Value = 1299.3572
$ws.Range("L58").Value = 6499.875
$ws.Range("M58").Value = -1096.3572
$ws.Range("N58").Value = -6905.875
$ws.Range("H99").Value = 4394.6875
$ws.Range("I99").Value = 2087.4285
$ws.Range("J99").Value = 6189.222
$ws.Range("K99").Value = 2087.4285
$ws.Range("L99").Value = 6189.222
$ws.Range("M99").Value = -589.4285
$ws.Range("N99").Value = -9185.222
$ws.Range("H126").Value = 4394.6875
$ws.Range("I126").Value = 2087.4285
$ws.Range("J126").Value = 6189.222
$ws.Range("K126").Value = 6262.2855
$ws.Range("L126").Value = 18567.666
$ws.Range("M126").Value = -3792.2855
$ws.Range("N126").Value = -23507.666
$ws.Range("H132").Value = 4120.6772
$ws.Range("I132").Value = 3329.9524
$ws.Range("J132").Value = 5781.2
$ws.Range("K132").Value = 9989.8572
$ws.Range("L132").Value = 17343.6
$ws.Range("M132").Value = -7459.8572
$ws.Range("N132").Value = -22403.6
$ws.Range("H134").Value = 1284.742
$ws.Range("I134").Value = 809.92
$ws.Range("J134").Value = 3263.1667
$ws.Range("K134").Value = 2429.76
$ws.Range("L134").Value = 9789.500100000001
$ws.Range("M134").Value = 105.2400000000002
$ws.Range("N134").Value = -14859.5001
$ws.Range("H136").Value = 2455.0278
$ws.Range("I136").Value = 1299.3572
$ws.Range("J136").Value = 6499.875
$ws.Range("K136").Value = 3898.0716
$ws.Range("L136").Value = 19499.625
$ws.Range("M136").Value = -1348.0716
$ws.Range("N136").Value = -24599.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 12082.35
$ws.Range("I34").Value = 23754
$ws.Range("J34").Value = 8191.8
$ws.Range("K34").Value = 71262
$ws.Range("L34").Value = 24575.4
$ws.Range("M34").Value = -71178
$ws.Range("N34").Value = -24743.4
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H49").Value = 4000
$ws.Range("J49").Value = 4000
$ws.Range("L49").Value = 12000
$ws.Range("N49").Value = -12312
$ws.Range("H64").Value = 1848.6364
$ws.Range("J64").Value = 2989.1667
$ws.Range("L64").Value = 8967.500100000001
$ws.Range("N64").Value = -9507.500100000001
$ws.Range("H67").Value = 1848.6364
$ws.Range("J67").Value = 2989.1667
$ws.Range("L67").Value = 8967.500100000001
$ws.Range("N67").Value = -10839.5001
$ws.Range("H68").Value = 1207.3383
$ws.Range("I68").Value = 1014.8182
$ws.Range("J68").Value = 1388.8572
$ws.Range("K68").Value = 3044.4546
$ws.Range("L68").Value = 4166.571599999999
$ws.Range("M68").Value = -2233.4546
$ws.Range("N68").Value = -5788.571599999999
$ws.Range("H71").Value = 1207.3383
$ws.Range("I71").Value = 1014.8182
$ws.Range("J71").Value = 1388.8572
$ws.Range("K71").Value = 9133.363800000001
$ws.Range("L71").Value = 12499.7148
$ws.Range("M71").Value = -5077.363800000001
$ws.Range("N71").Value = -20611.7148
$ws.Range("H94").Value = 999
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H131").Value = 785.72
$ws.Range("I131").Value = 476.30768
$ws.Range("J131").Value = 831.95404
$ws.Range("K131").Value = 1428.92304
$ws.Range("L131").Value = 2495.86212
$ws.Range("M131").Value = 3611.07696
$ws.Range("N131").Value = -12575.86212
$ws.Range("H133").Value = 8104.6665
$ws.Range("I133").Value = 9757
$ws.Range("K133").Value = 29271
$ws.Range("M133").Value = -24211

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2008
$ws.Range("I6").Value = 2008
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2008
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1895
$ws.Range("N6").ClearContents()
$ws.Range("H16").Value = 2008
$ws.Range("I16").Value = 2008
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2008
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1758
$ws.Range("N16").ClearContents()
$ws.Range("H132").Value = 2127.9092
$ws.Range("I132").Value = 1689.5306
$ws.Range("J132").Value = 5708
$ws.Range("K132").Value = 5068.5918
$ws.Range("L132").Value = 17124
$ws.Range("M132").Value = -2538.5918
$ws.Range("N132").Value = -22184
$ws.Range("H134").Value = 41640.668
$ws.Range("J134").Value = 41640.668
$ws.Range("L134").Value = 124922.004
$ws.Range("N134").Value = -129992.004

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2610.8914
$ws.Range("I136").Value = 1029.7
$ws.Range("J136").Value = 5575.625
$ws.Range("K136").Value = 3089.1
$ws.Range("L136").Value = 16726.875
$ws.Range("M136").Value = -539.1000000000004
$ws.Range("N136").Value = -21826.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51040
$ws.Range("H126").Value = 1333052.9
$ws.Range("I126").Value = 1825.2
$ws.Range("K126").Value = 5475.6
$ws.Range("M126").Value = -3005.6
$ws.Range("H132").Value = 4354.2
$ws.Range("I132").Value = 2580
$ws.Range("J132").Value = 4999.364
$ws.Range("K132").Value = 7740
$ws.Range("L132").Value = 14998.092
$ws.Range("M132").Value = -5210
$ws.Range("N132").Value = -20058.092
